# Refresh cryptos list (GitHub Actions scheduled update).
# Price/volume columns are stored as text (Price has thousands-separator
# dots, e.g. "29.211.40", so it can never be a real number) - values that
# look numeric are written with a leading "'" so Excel keeps them as text
# instead of silently re-typing the cell as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.211.40'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.842.59'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').Value = '''0.9999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''242.69'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').Value = '''0.6629'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''44.75'
$ws.Range('E8').Value = '  +6.35%  '
$ws.Range('D9').Value = '''0.07460'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = '''0.2958'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '''0.07747'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.840.68'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '''5.018'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '''0.6735'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '''83.41'
$ws.Range('E16').Value = '  -3.31%  '
$ws.Range('D17').Value = '''6.176'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('E18').Value = '  +5.18%  '
$ws.Range('D19').Value = '29.208.26'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '2.094.05'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = '''227.03'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '''7.194'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '''158.76'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').Value = '''8.628'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '''0.1395'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''1.209'
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''4.042'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').Value = '''0.05383'
$ws.Range('D35').Value = '''1.860'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').Value = '''0.7462'
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').Value = '''1.159'
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('D38').Value = '''2.652'
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('D39').Value = '1.300.61'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').Value = '''0.01798'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '''6.365'
$ws.Range('E42').Value = '  +6.89%  '
$ws.Range('D43').Value = '''0.9061'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('D44').Value = '''0.08331'
$ws.Range('E44').Value = '  +5.61%  '
$ws.Range('D45').Value = '''0.9999'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '''103.58'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = '1.991.44'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').Value = '''65.11'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '''1.752'
$ws.Range('E51').Value = '  -1.03%  '
